# Katalog guncellendi - Per 04.12.2025 12:02:47,87
#
# Applies the edit described in the commit:
#  - The product description for "Regular Fit Kot Pantolon" (rows 40-48)
#    and the two "Straight Kot Pantolon" rows (38-39) is shortened by
#    dropping the trailing "Ürünümüz tekli olarak satın alınabilir."
#    clause. This both replaces the old standalone description string
#    (previously only used by rows 38-39) and repoints rows 40-48 away
#    from the longer description (still used unchanged by rows 49-52)
#    to this new shared text.
#  - The sheet's active window scroll position / selection is moved from
#    C118:C124 to E38:E48 (matching where the edit was made).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newDescription = "%98 pamuk ve %2 spandex. 30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Belirtilen fiyatlar adet fiyatıdır."

$rows = 38..48
foreach ($r in $rows) {
    $ws.Range("E$r").Value2 = $newDescription
}

# Move the view / selection to reflect where the edit happened.
$ws.Range("E38:E48").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
